$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy cell formatting (style) for the A and E columns of the new rows 237-239
# from existing row 236, which already carries the correct styles (s="1" bold/border
# style for column A, and s="2" date style for column E).
$ws.Range('A236').Copy() | Out-Null
$ws.Range('A237:A239').PasteSpecial(-4122) | Out-Null
$ws.Range('E236').Copy() | Out-Null
$ws.Range('E237:E239').PasteSpecial(-4122) | Out-Null

# Row 234
$ws.Range('A234').Value = 232
$ws.Range('B234').Value = 6978437
$ws.Range('C234').Value = 'Bulgaria First League'
$ws.Range('D234').Value = 'Bulgaria First League'
$ws.Range('E234').Value = 45390.45833333334
$ws.Range('F234').Value = 'Botev Vratsa'
$ws.Range('G234').Value = 'Krumovgrad'
$ws.Range('H234').Value = 0
$ws.Range('I234').Value = 0
$ws.Range('J234').Value = 'D'
$ws.Range('K234').Value = 3
$ws.Range('L234').Value = 3.2
$ws.Range('M234').Value = 2.4
$ws.Range('N234').Value = 3.1
$ws.Range('O234').Value = 3.2
$ws.Range('P234').Value = 2.45
$ws.Range('Q234').Value = 0.25
$ws.Range('R234').Value = 1.775
$ws.Range('S234').Value = 2.1
$ws.Range('T234').Value = 2.25
$ws.Range('U234').Value = 2.05
$ws.Range('V234').Value = 1.8
$ws.Range('W234').Value = -1
$ws.Range('X234').Value = 2.2
$ws.Range('Y234').Value = -1
$ws.Range('Z234').Value = 0.3875
$ws.Range('AA234').Value = -0.5
$ws.Range('AB234').Value = -1
$ws.Range('AC234').Value = 0.8

# Row 235
$ws.Range('A235').Value = 233
$ws.Range('B235').Value = 6978448
$ws.Range('C235').Value = 'Bulgaria First League'
$ws.Range('D235').Value = 'Bulgaria First League'
$ws.Range('E235').Value = 45390.5625
$ws.Range('F235').Value = 'Lokomotiv Plovdiv'
$ws.Range('G235').Value = 'Cherno More Varna'
$ws.Range('H235').Value = 1
$ws.Range('I235').Value = 0
$ws.Range('J235').Value = 'H'
$ws.Range('K235').Value = 2.9
$ws.Range('L235').Value = 3.1
$ws.Range('M235').Value = 2.5
$ws.Range('N235').Value = 3.2
$ws.Range('O235').Value = 3
$ws.Range('P235').Value = 2.45
$ws.Range('Q235').Value = 0.25
$ws.Range('R235').Value = 1.775
$ws.Range('S235').Value = 2.1
$ws.Range('T235').Value = 2
$ws.Range('U235').Value = 1.8
$ws.Range('V235').Value = 2.05
$ws.Range('W235').Value = 2.2
$ws.Range('X235').Value = -1
$ws.Range('Y235').Value = -1
$ws.Range('Z235').Value = 0.7749999999999999
$ws.Range('AA235').Value = -1
$ws.Range('AB235').Value = -1
$ws.Range('AC235').Value = 1.05

# Row 236
$ws.Range('A236').Value = 234
$ws.Range('B236').Value = 6978449
$ws.Range('C236').Value = 'Bulgaria First League'
$ws.Range('D236').Value = 'Bulgaria First League'
$ws.Range('E236').Value = 45391.52083333334
$ws.Range('F236').Value = 'Lokomotiv 1929 Sofia'
$ws.Range('G236').Value = 'Beroe'
$ws.Range('K236').Value = 2.25
$ws.Range('L236').Value = 3
$ws.Range('M236').Value = 3.5
$ws.Range('N236').Value = 2.7
$ws.Range('O236').Value = 3
$ws.Range('P236').Value = 2.8
$ws.Range('Q236').Value = 0
$ws.Range('R236').Value = 1.85
$ws.Range('S236').Value = 2
$ws.Range('T236').Value = 2
$ws.Range('U236').Value = 1.8
$ws.Range('V236').Value = 2.05
$ws.Range('W236').Value = 0
$ws.Range('X236').Value = 0
$ws.Range('Y236').Value = 0
$ws.Range('Z236').Value = 0
$ws.Range('AA236').Value = 0

# Row 237
$ws.Range('A237').Value = 235
$ws.Range('B237').Value = 6978455
$ws.Range('C237').Value = 'Bulgaria First League'
$ws.Range('D237').Value = 'Bulgaria First League'
$ws.Range('E237').Value = 45394.59375
$ws.Range('F237').Value = 'CSKA Sofia'
$ws.Range('G237').Value = 'Pirin Blagoevgrad'
$ws.Range('K237').Value = 1.166
$ws.Range('L237').Value = 7
$ws.Range('M237').Value = 15
$ws.Range('N237').Value = 1.142
$ws.Range('O237').Value = 7
$ws.Range('P237').Value = 17
$ws.Range('Q237').Value = -2
$ws.Range('R237').Value = 1.925
$ws.Range('S237').Value = 1.925
$ws.Range('T237').Value = 2.75
$ws.Range('U237').Value = 1.875
$ws.Range('V237').Value = 1.975
$ws.Range('W237').Value = 0
$ws.Range('X237').Value = 0
$ws.Range('Y237').Value = 0
$ws.Range('Z237').Value = 0
$ws.Range('AA237').Value = 0

# Row 238
$ws.Range('A238').Value = 236
$ws.Range('B238').Value = 6978390
$ws.Range('C238').Value = 'Bulgaria First League'
$ws.Range('D238').Value = 'Bulgaria First League'
$ws.Range('E238').Value = 45395.48958333334
$ws.Range('F238').Value = 'Slavia Sofia'
$ws.Range('G238').Value = 'FC Hebar Pazardzhik'
$ws.Range('K238').Value = 1.25
$ws.Range('L238').Value = 6
$ws.Range('M238').Value = 10
$ws.Range('N238').Value = 1.363
$ws.Range('O238').Value = 5.25
$ws.Range('P238').Value = 7
$ws.Range('Q238').Value = -1.25
$ws.Range('R238').Value = 1.875
$ws.Range('S238').Value = 1.975
$ws.Range('T238').Value = 2.25
$ws.Range('U238').Value = 1.825
$ws.Range('V238').Value = 2.025
$ws.Range('W238').Value = 0
$ws.Range('X238').Value = 0
$ws.Range('Y238').Value = 0
$ws.Range('Z238').Value = 0
$ws.Range('AA238').Value = 0

# Row 239
$ws.Range('A239').Value = 237
$ws.Range('B239').Value = 6978457
$ws.Range('C239').Value = 'Bulgaria First League'
$ws.Range('D239').Value = 'Bulgaria First League'
$ws.Range('E239').Value = 45396.45833333334
$ws.Range('F239').Value = 'Etar 1924 Veliko Tarnovo'
$ws.Range('G239').Value = 'Krumovgrad'
$ws.Range('K239').Value = 4
$ws.Range('L239').Value = 3.2
$ws.Range('M239').Value = 2
$ws.Range('N239').Value = 3.5
$ws.Range('O239').Value = 3.1
$ws.Range('P239').Value = 2.2
$ws.Range('Q239').Value = 0.25
$ws.Range('R239').Value = 1.975
$ws.Range('S239').Value = 1.875
$ws.Range('T239').Value = 2
$ws.Range('U239').Value = 1.8
$ws.Range('V239').Value = 2.05
$ws.Range('W239').Value = 0
$ws.Range('X239').Value = 0
$ws.Range('Y239').Value = 0
$ws.Range('Z239').Value = 0
$ws.Range('AA239').Value = 0
